$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Bad Drivers" summary row (row 3) and Totals row (row 4)
$ws.Range("C3").Value = 295
$ws.Range("D3").Value = 96
$ws.Range("C4").Value = 295

# Update aggregated driver stats row (row 12)
$ws.Range("B12").Value = 59372
$ws.Range("E12").Value = 97
$ws.Range("F12").Value = 59433
